$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 217 (shifts old rows 217-241 down to 218-242)
$ws.Rows.Item(217).Insert()

# Populate the newly inserted row 217 with the weekly record
$ws.Cells.Item(217, 1).Value = 8
$ws.Cells.Item(217, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(217, 3).Value = 'Coquimbo'
$ws.Cells.Item(217, 4).Value = 45077
$ws.Cells.Item(217, 5).Value = 4
$ws.Cells.Item(217, 6).Value = 100112001
$ws.Cells.Item(217, 7).Value = 'Berenjena'
$ws.Cells.Item(217, 8).Value = 'Sin especificar'
$ws.Cells.Item(217, 9).Value = 'Primera'
$ws.Cells.Item(217, 10).Value = 400
$ws.Cells.Item(217, 11).Value = 8500
$ws.Cells.Item(217, 12).Value = 9000
$ws.Cells.Item(217, 13).Value = 8750
$ws.Cells.Item(217, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(217, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(217, 16).Value = 175
$ws.Cells.Item(217, 17).Value = 50
$ws.Cells.Item(217, 18).Value = 'Hortaliza'
